$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Open/Close/High/Low/Shares-Outstanding values per row (2..24), and the
# ticker string for column I is unified to "APP" for every data row (removing
# the need for the stray per-row ticker strings that had been injected into
# the shared-strings table from "extra files").
$rows = @(
    @{ Row = 2;  D = 70;                 E = 58.0099983215332;  F = 71.51000213623047;  G = 55.70500183105469;  H = 307636373 }
    @{ Row = 3;  D = 70;                 E = 58.0099983215332;  F = 71.51000213623047;  G = 55.70500183105469;  H = 307636373 }
    @{ Row = 4;  D = 70;                 E = 58.0099983215332;  F = 71.51000213623047;  G = 55.70500183105469;  H = 307636373 }
    @{ Row = 5;  D = 70;                 E = 58.0099983215332;  F = 71.51000213623047;  G = 55.70500183105469;  H = 307636373 }
    @{ Row = 6;  D = 70;                 E = 58.0099983215332;  F = 71.51000213623047;  G = 55.70500183105469;  H = 307636373 }
    @{ Row = 7;  D = 70;                 E = 58.0099983215332;  F = 71.51000213623047;  G = 55.70500183105469;  H = 307636373 }
    @{ Row = 8;  D = 74.27999877929688;  E = 61.47000122070312; F = 74.90000152587891;  G = 58.15000152587891;  H = 307636373 }
    @{ Row = 9;  D = 72.58999633789062;  E = 98.25;              F = 100.5;              G = 70.05000305175781;  H = 307636373 }
    @{ Row = 10; D = 95;                 E = 64.41999816894531; F = 98.13999938964844;  G = 57;                 H = 307636373 }
    @{ Row = 11; D = 55.2400016784668;   E = 38.15000152587891; F = 58.27000045776367;  G = 36.59999847412109;  H = 307636373 }
    @{ Row = 12; D = 34.68999862670898;  E = 35.54000091552734; F = 38.70000076293945;  G = 30.63100051879883;  H = 307636373 }
    @{ Row = 13; D = 19.6200008392334;   E = 16.95999908447266; F = 21.79000091552734;  G = 16.46599960327148;  H = 307636373 }
    @{ Row = 14; D = 10.8100004196167;   E = 12.69999980926514; F = 12.71000003814697;  G = 9.215000152587891;  H = 307636373 }
    @{ Row = 15; D = 15.61999988555908;  E = 17;                 F = 17.34000015258789;  G = 15.25;              H = 307636373 }
    @{ Row = 16; D = 25.64999961853028;  E = 31.39999961853028; F = 31.88999938964844;  G = 25.23999977111816;  H = 307636373 }
    @{ Row = 17; D = 40.15000152587891;  E = 36.43999862670898; F = 42.22000122070312;  G = 34.45000076293945;  H = 307636373 }
    @{ Row = 18; D = 39.40999984741211;  E = 41.13000106811523; F = 47.04000091552734;  G = 37.40000152587891;  H = 307636373 }
    @{ Row = 19; D = 69.04000091552734;  E = 70.56999969482422; F = 79.55000305175781;  G = 65.62000274658203;  H = 307636373 }
    @{ Row = 20; D = 85.12000274658203;  E = 77.09999847412109; F = 91.91000366210938;  G = 73.08499908447266;  H = 307636373 }
    @{ Row = 21; D = 130.7899932861328;  E = 169.3899993896484; F = 174.4600067138672;  G = 128.0099945068359;  H = 307636373 }
    @{ Row = 22; D = 331.0599975585937;  E = 369.5899963378906; F = 385.6600036621094;  G = 307.0679931640625;  H = 307636373 }
    @{ Row = 23; D = 264.2200012207031;  E = 269.3099975585937; F = 299.7200012207031;  G = 200.5;              H = 307636373 }
    @{ Row = 24; D = 347;                E = 390.7000122070313; F = 397.9200134277344;  G = 325.5799865722656;  H = 307636373 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = "APP"
}
